# Fix planning email, restore professional template, improve task deduplication
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: correct ticket_id (was a bogus phone-like number) and swap in the
#     real technician/patente/cliente values (was placeholder NANO data) ---
$ws.Cells.Item(2, 1).Value = 46011
$ws.Cells.Item(2, 2).Value = 7
$ws.Cells.Item(2, 9).Value = "Juan Perez"
$ws.Cells.Item(2, 10).Value = "HCCR38"
$ws.Cells.Item(2, 11).Value = "CIAL_ALIMENTOS"

# --- Row 3: same date/id fix, plus real technician/patente/cliente values
#     (was placeholder CHRIS data) ---
$ws.Cells.Item(3, 1).Value = 46011
$ws.Cells.Item(3, 2).Value = 8
$ws.Cells.Item(3, 9).Value = "Pedro Pascal"
$ws.Cells.Item(3, 10).Value = "BSBJ91"
$ws.Cells.Item(3, 11).Value = "CIAL_ALIMENTOS"

# --- Row 4: new ticket (previously blank) ---
$ws.Cells.Item(4, 1).Value = 46011
$ws.Cells.Item(4, 2).Value = 9
$ws.Cells.Item(4, 3).Value = "inmediata"
$ws.Cells.Item(4, 4).Value = "Soporte"
$ws.Cells.Item(4, 5).Value = "Botón Alámbrico Tablero"
$ws.Cells.Item(4, 6).Value = "RODRIGO DE ARAYA 2821, MACUL"
$ws.Cells.Item(4, 7).Value = "MACUL"
$ws.Cells.Item(4, 8).Value = "Región Metropolitana de Santiago."
$ws.Cells.Item(4, 9).Value = "Juan Perez"
$ws.Cells.Item(4, 10).Value = "VDVP13"
$ws.Cells.Item(4, 11).Value = "COMERCIAL_DIBOR"

# --- Row 5: new ticket (previously blank) ---
$ws.Cells.Item(5, 1).Value = 46011
$ws.Cells.Item(5, 2).Value = 10
$ws.Cells.Item(5, 3).Value = "alta"
$ws.Cells.Item(5, 4).Value = "Soporte"
$ws.Cells.Item(5, 5).Value = "MDVR 4ch sin IA"
$ws.Cells.Item(5, 6).Value = "EL VENTISQUERO 1250, RENCA"
$ws.Cells.Item(5, 7).Value = "RENCA"
$ws.Cells.Item(5, 8).Value = "Región Metropolitana de Santiago."
$ws.Cells.Item(5, 9).Value = "Juan Perez"
$ws.Cells.Item(5, 10).Value = "KHSJ34"
$ws.Cells.Item(5, 11).Value = "SOPROLE"

# --- Row 6: new ticket (previously blank) ---
$ws.Cells.Item(6, 1).Value = 46011
$ws.Cells.Item(6, 2).Value = 11
$ws.Cells.Item(6, 3).Value = "normal"
$ws.Cells.Item(6, 4).Value = "Instalación"
$ws.Cells.Item(6, 5).Value = "Antena GPS"
$ws.Cells.Item(6, 6).Value = "AV. LO ESPEJO 1300, MAIPU"
$ws.Cells.Item(6, 7).Value = "MAIPU"
$ws.Cells.Item(6, 8).Value = "Región Metropolitana de Santiago."
$ws.Cells.Item(6, 9).Value = "Pedro Pascal"
$ws.Cells.Item(6, 10).Value = "GENERADOR_SANTIAGO_19"
$ws.Cells.Item(6, 11).Value = "LUREYE"

# --- Row 7: new ticket (previously blank), deduplicated against row 6
#     (same client/address, distinct generator asset tag) ---
$ws.Cells.Item(7, 1).Value = 46011
$ws.Cells.Item(7, 2).Value = 12
$ws.Cells.Item(7, 3).Value = "normal"
$ws.Cells.Item(7, 4).Value = "Instalación"
$ws.Cells.Item(7, 5).Value = "Antena GPS"
$ws.Cells.Item(7, 6).Value = "AV. LO ESPEJO 1300, MAIPU"
$ws.Cells.Item(7, 7).Value = "MAIPU"
$ws.Cells.Item(7, 8).Value = "Región Metropolitana de Santiago."
$ws.Cells.Item(7, 9).Value = "Pedro Pascal"
$ws.Cells.Item(7, 10).Value = "GENERADOR_SANTIAGO_18"
$ws.Cells.Item(7, 11).Value = "LUREYE"

# Move/leave the cursor where the author left it after finishing data entry.
$ws.Range("B8").Select() | Out-Null
